# Update "想去人数" (F2:F5) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 460
    $ws.Range("F3").Value = 3278
    $ws.Range("F4").Value = 84
    $ws.Range("F5").Value = 653
}
